# Apply the edit described by the diff:
#  - Append a new data row (row 11) to Hoja1:
#       A11 = "Salesforce", B11 = "Layout", C11 = "Objeto__c", D11 = "Página Objeto"
#  - This introduces three new shared strings: "Layout", "Página Objeto", "Objeto__c"
#  - Selection/active cell moves to C13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A11").Value = "Salesforce"
$ws.Range("B11").Value = "Layout"
$ws.Range("D11").Value = "Página Objeto"
$ws.Range("C11").Value = "Objeto__c"

$ws.Range("C13").Select()
